$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the H1 title.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# 2. Replace the final (italic, "Create a feature image..." prompt) paragraph
#    with a new bold "Play Glitz Free Slot Game | Simple Gameplay Design"
#    paragraph followed by an italic paragraph holding the (re-worded) meta
#    description text.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Create a feature image for Glitz*") {
        $p.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Glitz Free Slot Game | Simple Gameplay Design</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Glitz slot game. Play for free and experience simple gameplay with 60 paylines, free spin bonuses, and a relaxing musical choice.</w:t></w:r></w:p>") | Out-Null
        break
    }
}
